# Adds a new column "Numero de Concepto" right before the existing "Path"
# column. The new column's value, for each data row, is a running count of
# how many times that row's "Path" value has been seen so far (1-based),
# i.e. an occurrence counter per distinct PDF path, in row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 183
$pathCol = 13   # column M ("Path") before the insert

# 1) Read existing Path values (column M) before shifting anything.
$paths = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $paths[$r] = $ws.Cells.Item($r, $pathCol).Value2
}

# 2) Insert a new blank column at M; old M ("Path") becomes N.
$ws.Columns("M").Insert()

# 3) Header for the new column.
$ws.Range("M1").Value = "Numero de Concepto"

# 4) Fill the new column with a running per-path occurrence counter.
$counts = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $path = $paths[$r]
    if ($counts.ContainsKey($path)) {
        $counts[$path] = $counts[$path] + 1
    } else {
        $counts[$path] = 1
    }
    $ws.Cells.Item($r, $pathCol).Value = [string]$counts[$path]
}
